$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$shape = $s.Shapes.Item(1)
$shape.TextFrame.TextRange.Text = "Hệ thống theo dõi và xử lý tác vụ của" + [char]13 + "Entity Framework"
